$d = $word.ActiveDocument

# Change 1: placeholder key changed from "#OverviewTable_..." to "#Country_..."
# (was split across two runs: " #Overview" and "Table")
$d.Content.Find.Execute("#OverviewTable", $true, $false, $false, $false, $false,
                         $true, 1, $false, "#Country", 2)

# Change 2: the run containing "sit" was wrapped with proofErr gramStart/gramEnd
# grammar-check markers, splitting "...Lorem ipsum dolor " / "sit" / " amet. Lorem..."
# into three runs. Re-typing the whole sentence as one pass removes the stale
# proofErr bookmarks and collapses the text back into a single run.
$d.Content.Find.Execute("Stet clita kasd gubergren, no sea takimata sanctus est Lorem ipsum dolor sit amet.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Stet clita kasd gubergren, no sea takimata sanctus est Lorem ipsum dolor sit amet.", 2)
